$d = $word.ActiveDocument

function Insert-Xml($range, $bodyFragment) {
    $pkg = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + $bodyFragment + '</w:document></pkg:xmlData></pkg:part></pkg:package>'
    $range.InsertXML($pkg)
}

# ---------------------------------------------------------------------
# 1. Insert a centered, bold 20pt "CRC Cards" title paragraph followed
#    by a blank paragraph at the very start of the document.
# ---------------------------------------------------------------------
$titleFrag = '<w:body><w:p><w:pPr><w:jc w:val="center"/><w:rPr><w:b/><w:bCs/><w:sz w:val="40"/><w:szCs w:val="40"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="40"/><w:szCs w:val="40"/></w:rPr><w:t>CRC Cards</w:t></w:r></w:p><w:p/></w:body>'
$titleRange = $d.Range(0, 0)
Insert-Xml $titleRange $titleFrag

# ---------------------------------------------------------------------
# 2. Mark the last cursor position (the empty cell beside "Tax.java")
#    with a "_GoBack" bookmark, matching Word's own auto-bookmark.
# ---------------------------------------------------------------------
$firstTable = $d.Tables(1)
$goBackCell = $firstTable.Cell(1, 2)
$d.Bookmarks.Add("_GoBack", $goBackCell.Range)

# ---------------------------------------------------------------------
# 3. Append three new CRC-card tables (GUI, Property, Property
#    management IMPL) after the existing last table.
# ---------------------------------------------------------------------
$newTablesXml = @'
<w:p/><w:p/><w:tbl><w:tblPr><w:tblStyle w:val="TableGrid"/><w:tblW w:w="0" w:type="auto"/><w:tblLook w:val="04A0" w:firstRow="1" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:noHBand="0" w:noVBand="1"/></w:tblPr><w:tblGrid><w:gridCol w:w="4505"/><w:gridCol w:w="4505"/></w:tblGrid><w:tr><w:tc><w:tcPr><w:tcW w:w="4508" w:type="dxa"/><w:tcBorders><w:top w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:left w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:bottom w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:right w:val="single" w:sz="4" w:space="0" w:color="auto"/></w:tcBorders><w:hideMark/></w:tcPr><w:p><w:r><w:t>GUI</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="4508" w:type="dxa"/><w:tcBorders><w:top w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:left w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:bottom w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:right w:val="single" w:sz="4" w:space="0" w:color="auto"/></w:tcBorders><w:hideMark/></w:tcPr><w:p><w:r><w:t xml:space="preserve">Collaboration </w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="4508" w:type="dxa"/><w:tcBorders><w:top w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:left w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:bottom w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:right w:val="single" w:sz="4" w:space="0" w:color="auto"/></w:tcBorders></w:tcPr><w:p><w:r><w:t>Creates a GUI,</w:t></w:r></w:p><w:p><w:r><w:t>That user can input data into,</w:t></w:r></w:p><w:p><w:r><w:t>I.e., their address, Eircode, name, house, price etc.</w:t></w:r></w:p><w:p/></w:tc><w:tc><w:tcPr><w:tcW w:w="4508" w:type="dxa"/><w:tcBorders><w:top w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:left w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:bottom w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:right w:val="single" w:sz="4" w:space="0" w:color="auto"/></w:tcBorders><w:hideMark/></w:tcPr><w:p><w:r><w:t>Test</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Property management </w:t></w:r></w:p><w:p><w:r><w:t>Tax</w:t></w:r></w:p></w:tc></w:tr></w:tbl><w:p><w:pPr><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr></w:p><w:p/><w:p/><w:p/><w:tbl><w:tblPr><w:tblStyle w:val="TableGrid"/><w:tblW w:w="0" w:type="auto"/><w:tblLook w:val="04A0" w:firstRow="1" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:noHBand="0" w:noVBand="1"/></w:tblPr><w:tblGrid><w:gridCol w:w="4505"/><w:gridCol w:w="4505"/></w:tblGrid><w:tr><w:tc><w:tcPr><w:tcW w:w="4508" w:type="dxa"/><w:tcBorders><w:top w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:left w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:bottom w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:right w:val="single" w:sz="4" w:space="0" w:color="auto"/></w:tcBorders><w:hideMark/></w:tcPr><w:p><w:r><w:t xml:space="preserve">Property </w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="4508" w:type="dxa"/><w:tcBorders><w:top w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:left w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:bottom w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:right w:val="single" w:sz="4" w:space="0" w:color="auto"/></w:tcBorders><w:hideMark/></w:tcPr><w:p><w:r><w:t xml:space="preserve">Collaboration </w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="4508" w:type="dxa"/><w:tcBorders><w:top w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:left w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:bottom w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:right w:val="single" w:sz="4" w:space="0" w:color="auto"/></w:tcBorders></w:tcPr><w:p><w:r><w:t>Creates new array list &lt;tax&gt;</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Sets owners </w:t></w:r></w:p><w:p><w:r><w:t>Sets address</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Sets Eircode </w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Sets value </w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Sets string, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>ppr</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> and location category</w:t></w:r></w:p><w:p/></w:tc><w:tc><w:tcPr><w:tcW w:w="4508" w:type="dxa"/><w:tcBorders><w:top w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:left w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:bottom w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:right w:val="single" w:sz="4" w:space="0" w:color="auto"/></w:tcBorders><w:hideMark/></w:tcPr><w:p><w:r><w:t>GUI</w:t></w:r></w:p><w:p><w:r><w:t>Tax</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Property management </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>impl</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p><w:r><w:t>Command line</w:t></w:r></w:p></w:tc></w:tr></w:tbl><w:p><w:pPr><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr></w:p><w:p/><w:p/><w:p/><w:tbl><w:tblPr><w:tblStyle w:val="TableGrid"/><w:tblW w:w="0" w:type="auto"/><w:tblLook w:val="04A0" w:firstRow="1" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:noHBand="0" w:noVBand="1"/></w:tblPr><w:tblGrid><w:gridCol w:w="4505"/><w:gridCol w:w="4505"/></w:tblGrid><w:tr><w:tc><w:tcPr><w:tcW w:w="4508" w:type="dxa"/><w:tcBorders><w:top w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:left w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:bottom w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:right w:val="single" w:sz="4" w:space="0" w:color="auto"/></w:tcBorders><w:hideMark/></w:tcPr><w:p><w:r><w:t>Property management IMPL</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="4508" w:type="dxa"/><w:tcBorders><w:top w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:left w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:bottom w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:right w:val="single" w:sz="4" w:space="0" w:color="auto"/></w:tcBorders><w:hideMark/></w:tcPr><w:p><w:r><w:t xml:space="preserve">Collaboration </w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="4508" w:type="dxa"/><w:tcBorders><w:top w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:left w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:bottom w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:right w:val="single" w:sz="4" w:space="0" w:color="auto"/></w:tcBorders></w:tcPr><w:p><w:r><w:t>This implements property management interface class.</w:t></w:r></w:p><w:p><w:r><w:t>Adds owner with name, address,</w:t></w:r></w:p><w:p><w:r><w:t>Value, and location category.</w:t></w:r></w:p><w:p/></w:tc><w:tc><w:tcPr><w:tcW w:w="4508" w:type="dxa"/><w:tcBorders><w:top w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:left w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:bottom w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:right w:val="single" w:sz="4" w:space="0" w:color="auto"/></w:tcBorders><w:hideMark/></w:tcPr><w:p><w:r><w:t xml:space="preserve">Property management interface </w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Property </w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Command line </w:t></w:r></w:p><w:p><w:r><w:t>GUI</w:t></w:r></w:p></w:tc></w:tr></w:tbl>
'@

$lastTable = $d.Tables($d.Tables.Count)
$insertPoint = $d.Range($lastTable.Range.End, $lastTable.Range.End)
$newTablesFrag = '<w:body>' + $newTablesXml + '<w:p/></w:body>'
Insert-Xml $insertPoint $newTablesFrag

# ---------------------------------------------------------------------
# 4. Remove the header text "CRC Cards- Leon Woods 19251727" but keep
#    the (now empty) paragraph that held it.
# ---------------------------------------------------------------------
$sec = $d.Sections(1)
$hdr = $sec.Headers(1)
$hdrPara = $hdr.Range.Paragraphs(1)
$hdrPara.Range.Text = ""

Write-Output "done"
